$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Summary"
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").Value = 0.3558718861209965
$wsSummary.Range("C2").Value = 0.06510416666666667
$wsSummary.Range("D2").Value = 0.8928571428571429
$wsSummary.Range("E2").Value = 0.1213592233009709
$wsSummary.Range("F2").Value = 0.2520161290322581
$wsSummary.Range("G2").Value = 0.5996309963099631
$wsSummary.Range("H2").Value = 0.8013643659711074
$wsSummary.Range("I2").Value = 25
$wsSummary.Range("J2").Value = 359
$wsSummary.Range("K2").Value = 175
$wsSummary.Range("L2").Value = 3

# ---------------------------------------------------------------------------
# Sheet 2: "Classification Report"
# ---------------------------------------------------------------------------
$wsReport = $wb.Worksheets.Item("Classification Report")

$wsReport.Range("B2").Value = 0.9831460674157303
$wsReport.Range("C2").Value = 0.3277153558052435
$wsReport.Range("D2").Value = 0.4915730337078651

$wsReport.Range("B3").Value = 0.06510416666666667
$wsReport.Range("C3").Value = 0.8928571428571429
$wsReport.Range("D3").Value = 0.1213592233009709

$wsReport.Range("B4").Value = 0.3558718861209965
$wsReport.Range("C4").Value = 0.3558718861209965
$wsReport.Range("D4").Value = 0.3558718861209965
$wsReport.Range("E4").Value = 0.3558718861209965

$wsReport.Range("B5").Value = 0.5241251170411985
$wsReport.Range("C5").Value = 0.6102862493311931
$wsReport.Range("D5").Value = 0.306466128504418

$wsReport.Range("B6").Value = 0.937407325029656
$wsReport.Range("C6").Value = 0.3558718861209965
$wsReport.Range("D6").Value = 0.4731282175310092

# ---------------------------------------------------------------------------
# Sheet 3: "Confusion Matrix"
# ---------------------------------------------------------------------------
$wsConfusion = $wb.Worksheets.Item("Confusion Matrix")

$wsConfusion.Range("B2").Value = 175
$wsConfusion.Range("C2").Value = 359
$wsConfusion.Range("B3").Value = 3
$wsConfusion.Range("C3").Value = 25
